$d = $word.ActiveDocument

# --- Step 1: Paragraph 2 ("-, Hiển thị câu hỏi lấy ra từ db") ---
# Replaced with new wording (split across 4 runs) AND a brand-new paragraph
# ("-, Mỗi lần trả lời...") is inserted right after it. We replace the whole
# paragraph range (including its trailing paragraph mark) with two <w:p>
# elements so the split happens in one atomic, order-preserving operation.
$p2 = $d.Paragraphs.Item(2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End)
$xml1 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">-, </w:t></w:r><w:r><w:t>Random lấy ra các câu hỏi trong db</w:t></w:r><w:r><w:t>(Tối đa 10 câu)</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">-, Mỗi lần trả lời sẽ kiểm tra kết quả và trả về kết quả đúng. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml1)

# --- Step 2: Paragraph "-, Xây dựng db" (now paragraph 6) gains two more
# runs (a long run of spaces starting with ":" and another run of spaces)
# right before the inline drawing. Inserting exactly at the shape boundary
# is unreliable, so we target the final character of the existing text
# ("b") and replace it with itself plus the two new runs - this keeps the
# new runs correctly ordered before the drawing.
$shp = $d.InlineShapes.Item(1)
$rLastChar = $d.Range($shp.Range.Start - 1, $shp.Range.Start)
$xml2 = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>b</w:t></w:r><w:r><w:t xml:space="preserve">:                                                                                                                                    </w:t></w:r><w:r><w:t xml:space="preserve">                                                                                                                           </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rLastChar.InsertXML($xml2)
